# "Generate Report for Handoff"
#
# b.md has now been handed off for zh-cn and de-de, so its Status moves
# from "Handed back: in sync with en-US" to "Ready for handoff", a new
# handoff file/time is recorded, and the Overview sheet rolls up the
# latest handoff date for b.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) now reflects the new handoff.
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"
$ovw.Range("D3").Value = "2016-24-18 10:24:57"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3) gets a new handoff file + datetime and
# its status flips to "Ready for handoff". The hyperlink display text
# for the "Latest Handoff File" cell needs updating too, so the whole
# hyperlink set on this sheet is rebuilt (the underlying URLs are kept
# identical to the originals).
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-18 10:24:54"

$zh.Hyperlinks.Delete()
$zhHl = $zh.Hyperlinks
$zhHl.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/a.md", "", "", "a.md")
$zhHl.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/a.md", "", "", ".md")
$zhHl.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/820165064aa34e84e813693a9a6ec975ebdef250/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhHl.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/642f3c6f3b3bed761ff11b5364b3bc48e575cae9/e2e/a.md", "", "", "a.md")
$zhHl.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27622e26f2034c9a01e060442264193e637420e6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhHl.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/b.md", "", "", "b.md")
$zhHl.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/b.md", "", "", ".md")
$zhHl.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/820165064aa34e84e813693a9a6ec975ebdef250/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhHl.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/642f3c6f3b3bed761ff11b5364b3bc48e575cae9/e2e/a.md", "", "", "a.md")
$zhHl.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27622e26f2034c9a01e060442264193e637420e6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, with de-de specific URLs.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-18 10:24:57"

$de.Hyperlinks.Delete()
$deHl = $de.Hyperlinks
$deHl.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/a.md", "", "", "a.md")
$deHl.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/a.md", "", "", ".md")
$deHl.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37e9cff96110580122ba9d1886e15de7b93d12f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$deHl.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b4c09958110b73bb3f8cf7f9c234a688976b5bbc/e2e/a.md", "", "", "a.md")
$deHl.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75b7602298e6c73c052bd6fdc98b5eb4104fa308/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$deHl.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/b.md", "", "", "b.md")
$deHl.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e/b.md", "", "", ".md")
$deHl.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37e9cff96110580122ba9d1886e15de7b93d12f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$deHl.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b4c09958110b73bb3f8cf7f9c234a688976b5bbc/e2e/a.md", "", "", "a.md")
$deHl.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75b7602298e6c73c052bd6fdc98b5eb4104fa308/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
